# Edit script for 5045867108_VCHAINS.xlsx
# This applies the changes described in the diff:
#  - Update TimeStamp (B1) and Company Id (B4)
#  - Update the 4 supplier/customer data rows (7-10) with new company data
#  - Delete the old row 11 (5th data row no longer present)
#  - Shrink the Data / DataConfidenceScore / DataRevenue named ranges accordingly
#  - Best-effort nudge of column widths B and F (cosmetic autofit side effect)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("B1").Value = 44598.6590277778      # TimeStamp
$ws.Range("B3").Value = "Chemours Co"         # Company Name
$ws.Range("B4").Value = 5045867108            # Company Id (matches file name)

# --- Row 7 (Apple Inc / Customer) only the "Days Since Last Update" changed
$ws.Range("I7").Value = 1466

# --- Row 8: Roku Inc (Supplier) -> Energy Fuels Inc (Customer) -------
$ws.Range("A8").Value = 4295862622
$ws.Range("B8").Value = "Energy Fuels Inc"
$ws.Range("C8").Value = "Public"
$ws.Range("D8").Value = "Customer"
$ws.Range("E8").Value = "United States of America"
$ws.Range("F8").Value = "Uranium"
$ws.Range("G8").Value = 0.666094594184154
$ws.Range("H8").Value = 44307
$ws.Range("I8").Value = 291
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 1658000
$ws.Range("M8").Value = 7
$ws.Range("N8").ClearContents()

# --- Row 9: Cypress Semiconductor Corp -> Pricewaterhousecoopers LLP -
$ws.Range("A9").Value = 5000017069
$ws.Range("B9").Value = "Pricewaterhousecoopers LLP"
$ws.Range("C9").Value = "Private"
$ws.Range("D9").Value = "Supplier"
$ws.Range("E9").Value = "United States of America"
$ws.Range("F9").Value = "Business Support Services"
$ws.Range("G9").Value = 0.31675424
$ws.Range("H9").Value = 43175
$ws.Range("I9").Value = 1423
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 6794563060
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

# --- Row 10: Optomec Inc -> DuPont Group Inc --------------------------
$ws.Range("A10").Value = 5037205994
$ws.Range("B10").Value = "DuPont Group Inc"
$ws.Range("C10").Value = "Private"
$ws.Range("D10").Value = "Supplier"
$ws.Range("E10").Value = "United States of America"
$ws.Range("F10").Value = "Commodity Chemicals"
$ws.Range("G10").Value = 0.2955984
$ws.Range("H10").Value = 42186
$ws.Range("I10").Value = 2412
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 2
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

# --- Remove the old 5th data row (Kioxia Holdings Corp) ---------------
$ws.Rows.Item(11).Delete()

# --- Shrink named ranges to match the smaller data block --------------
$wb.Names.Item("Data").RefersTo = "='Value Chains'!`$A`$7:`$N`$10"
$wb.Names.Item("DataConfidenceScore").RefersTo = "='Value Chains'!`$G`$7:`$G`$10"
$wb.Names.Item("DataRevenue").RefersTo = "='Value Chains'!`$I`$7:`$I`$10"

# --- Cosmetic: nudge column widths (closest achievable on this grid) --
$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(6).ColumnWidth = 24.83
